$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.303725242614746
$ws.Range("B1").Value = 1.927546262741089
$ws.Range("C1").Value = 3.892374753952026
$ws.Range("D1").Value = 0.892136812210083
$ws.Range("E1").Value = 0.7888664603233337
